# SCD0174 - Added "Goto SubNavBar Using Text Function"
#
# The test-scenario sheet documents a UI test step that navigates to a
# sub-navbar item by its text. The row-3 "SUB_NAVBAR" sample value ("Free")
# is relocated from column M (SUB_NAVBAR) to column O (TEXT1), and a new
# column P (TEXT2) value "Prospek" is added - together describing a
# "Goto ... Using Text" style lookup over two text fields. The RUN flag in
# A2 is also cleared (left blank, formatting kept), and the active
# selection/navigation moves to M3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the RUN indicator in A2 (keep its existing cell style/formatting).
$ws.Range("A2").ClearContents()

# Relocate the SUB_NAVBAR sample text ("Free") from M3 into O3 (TEXT1),
# and record the new sub-navbar target ("Prospek") into P3 (TEXT2).
$subNavBarValue = $ws.Range("M3").Value2
$ws.Range("M3").ClearContents()
$ws.Range("O3").Value = $subNavBarValue
$ws.Range("P3").Value = "Prospek"

# Reflect the "Goto ... Using Text" navigation by moving/selecting M3.
$excel.Goto($ws.Range("M3"), $true)
